$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Convention change to support multi-axle vehicles: add a "Truck_Amandla"
# sheet (between Bus_Makhulu and Trailer_Thwala) and a "Trailer_Kumanzi"
# sheet (after Trailer_Thwala). Both are built from the same 8-row template
# used by Bus_Makhulu, so we clone that sheet and adjust its data.
# ---------------------------------------------------------------------------

# --- Truck_Amandla: cloned from Bus_Makhulu, placed right after it ---------
$busTemplate = $wb.Worksheets.Item("Bus_Makhulu")
$busTemplate.Copy($null, $busTemplate)
$truck = $wb.Worksheets.Item("Bus_Makhulu (2)")
$truck.Name = "Truck_Amandla"

$truck.Range("H3").Value = "CAD_Truck_Amandla"
$truck.Range("H4").Value = "CAD_Truck_Amandla"
$truck.Range("F7").Value = 0.6
$truck.Range("G7").Value = 0.8
$truck.Range("H7").Value = 1
$truck.Range("H8").Value = 1

$truck.Activate()
$truck.Range("G23").Select()

# --- Trailer_Kumanzi: cloned from Bus_Makhulu too, placed after Trailer_Thwala ---
$trailerThwala = $wb.Worksheets.Item("Trailer_Thwala")
$busTemplate.Copy($null, $trailerThwala)
$kumanzi = $wb.Worksheets.Item("Bus_Makhulu (2)")
$kumanzi.Name = "Trailer_Kumanzi"

$kumanzi.Range("H3").Value = "CAD_Trailer_Kumanzi"
$kumanzi.Range("H4").Value = "CAD_Trailer_Kumanzi"
$kumanzi.Range("F7").Value = 1
$kumanzi.Range("G7").Value = 0.75
$kumanzi.Range("H7").Value = 0.055
$kumanzi.Range("H8").Value = 0.5

$kumanzi.Activate()
$kumanzi.Range("H8").Select()
